$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 25665.334
$ws.Range("J3").Value = 25665.334
$ws.Range("L3").Value = 25665.334
$ws.Range("N3").Value = -25893.334
$ws.Range("H43").Value = 2500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2500
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2500
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2638
$ws.Range("H69").Value = 15998.214
$ws.Range("I69").Value = 9993.5
$ws.Range("K69").Value = 29980.5
$ws.Range("M69").Value = -29106.5
$ws.Range("H72").Value = 15998.214
$ws.Range("I72").Value = 9993.5
$ws.Range("K72").Value = 89941.5
$ws.Range("M72").Value = -85573.5
$ws.Range("H98").Value = 1522.7858
$ws.Range("I98").Value = 1522.7858
$ws.Range("K98").Value = 1522.7858
$ws.Range("M98").Value = -24.78580000000011
$ws.Range("H102").Value = 25665.334
$ws.Range("J102").Value = 25665.334
$ws.Range("L102").Value = 25665.334
$ws.Range("N102").Value = -32155.334
$ws.Range("H122").Value = 1522.7858
$ws.Range("I122").Value = 1522.7858
$ws.Range("K122").Value = 4568.357400000001
$ws.Range("M122").Value = -2118.357400000001
$ws.Range("H132").Value = 2648.7942
$ws.Range("I132").Value = 2731.742
$ws.Range("K132").Value = 8195.226000000001
$ws.Range("M132").Value = -5665.226000000001
$ws.Range("H137").Value = 2181.55
$ws.Range("I137").Value = 1562.0869
$ws.Range("J137").Value = 3019.647
$ws.Range("K137").Value = 4686.2607
$ws.Range("L137").Value = 9058.940999999999
$ws.Range("M137").Value = -2136.2607
$ws.Range("N137").Value = -14158.941
$ws.Range("H138").Value = 4004.6438
$ws.Range("I138").Value = 3556.3076
$ws.Range("J138").Value = 4101.783
$ws.Range("K138").Value = 10668.9228
$ws.Range("L138").Value = 12305.349
$ws.Range("M138").Value = -5528.9228
$ws.Range("N138").Value = -22585.349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 3698.6
$ws.Range("I21").Value = 3500
$ws.Range("J21").Value = 3996.5
$ws.Range("K21").Value = 3500
$ws.Range("L21").Value = 3996.5
$ws.Range("M21").Value = -3126
$ws.Range("N21").Value = -4744.5
$ws.Range("H61").Value = 32261592
$ws.Range("I61").Value = 35717668
$ws.Range("J61").Value = 4876
$ws.Range("K61").Value = 35717668
$ws.Range("L61").Value = 4876
$ws.Range("M61").Value = -35717456
$ws.Range("N61").Value = -5300
$ws.Range("H74").Value = 25004138
$ws.Range("I74").Value = 26320020
$ws.Range("J74").Value = 2387
$ws.Range("K74").Value = 26320020
$ws.Range("L74").Value = 2387
$ws.Range("M74").Value = -26319146
$ws.Range("N74").Value = -4135
$ws.Range("H76").Value = 9666.666999999999
$ws.Range("J76").Value = 9666.666999999999
$ws.Range("L76").Value = 9666.666999999999
$ws.Range("N76").Value = -10342.667
$ws.Range("H77").Value = 25004138
$ws.Range("I77").Value = 26320020
$ws.Range("J77").Value = 2387
$ws.Range("K77").Value = 131600100
$ws.Range("L77").Value = 11935
$ws.Range("M77").Value = -131595732
$ws.Range("N77").Value = -20671
$ws.Range("H79").Value = 9666.666999999999
$ws.Range("J79").Value = 9666.666999999999
$ws.Range("L79").Value = 9666.666999999999
$ws.Range("N79").Value = -12006.667
$ws.Range("H102").Value = 1509.125
$ws.Range("I102").Value = 1299.7693
$ws.Range("K102").Value = 1299.7693
$ws.Range("M102").Value = 322.2307000000001
$ws.Range("H104").Value = 1974089.6
$ws.Range("J104").Value = 1974089.6
$ws.Range("L104").Value = 1974089.6
$ws.Range("N104").Value = -1981077.6
$ws.Range("H110").Value = 49640.953
$ws.Range("I110").Value = 60466.707
$ws.Range("J110").Value = 3631.5
$ws.Range("K110").Value = 60466.707
$ws.Range("L110").Value = 3631.5
$ws.Range("M110").Value = -58421.707
$ws.Range("N110").Value = -7721.5
$ws.Range("H136").Value = 32261592
$ws.Range("I136").Value = 35717668
$ws.Range("J136").Value = 4876
$ws.Range("K136").Value = 107153004
$ws.Range("L136").Value = 14628
$ws.Range("M136").Value = -107150454
$ws.Range("N136").Value = -19728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 37610.832
$ws.Range("J100").Value = 37610.832
$ws.Range("L100").Value = 37610.832
$ws.Range("N100").Value = -39774.832
$ws.Range("H107").Value = 30391.766
$ws.Range("I107").Value = 1017.8333
$ws.Range("K107").Value = 1017.8333
$ws.Range("M107").Value = 902.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 23810798
$ws.Range("J132").Value = 1221.25
$ws.Range("L132").Value = 3663.75
$ws.Range("N132").Value = -8723.75
$ws.Range("H141").Value = 241583.33
$ws.Range("J141").Value = 432216.66
$ws.Range("L141").Value = 432216.66
$ws.Range("N141").Value = -442576.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 12699.8
$ws.Range("J93").Value = 12699.8
$ws.Range("L93").Value = 38099.39999999999
$ws.Range("N93").Value = -41843.39999999999
$ws.Range("H97").Value = 2144
$ws.Range("I97").Value = 487.33334
$ws.Range("J97").Value = 2526.3076
$ws.Range("K97").Value = 1462.00002
$ws.Range("L97").Value = 7578.9228
$ws.Range("M97").Value = -966.0000199999999
$ws.Range("N97").Value = -8570.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8287.223
$ws.Range("J92").Value = 8287.223
$ws.Range("L92").Value = 8287.223
$ws.Range("N92").Value = -12031.223
$ws.Range("H97").Value = 981.7273
$ws.Range("I97").Value = 743.2857
$ws.Range("K97").Value = 743.2857
$ws.Range("M97").Value = -247.2857
$ws.Range("H132").Value = 4314942.5
$ws.Range("I132").Value = 4812075
$ws.Range("J132").Value = 6465
$ws.Range("K132").Value = 14436225
$ws.Range("L132").Value = 19395
$ws.Range("M132").Value = -14433695
$ws.Range("N132").Value = -24455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 883.55554
$ws.Range("I81").Value = 844
$ws.Range("K81").Value = 1688
$ws.Range("M81").Value = -627
$ws.Range("H84").Value = 883.55554
$ws.Range("I84").Value = 844
$ws.Range("K84").Value = 8440
$ws.Range("M84").Value = -3136
$ws.Range("H107").Value = 2430.7144
$ws.Range("I107").Value = 1412.2
$ws.Range("J107").Value = 2996.5557
$ws.Range("K107").Value = 4236.6
$ws.Range("L107").Value = 8989.667099999999
$ws.Range("M107").Value = -2316.6
$ws.Range("N107").Value = -12829.6671
$ws.Range("H136").Value = 18520500
$ws.Range("I136").Value = 20834772
$ws.Range("K136").Value = 62504316
$ws.Range("M136").Value = -62501766
